$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style used by the other
# header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I (I0) and J (IF), rows 2-35.
$data = @(
    @(5, 5),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(1, 2),
    @(7, 7),
    @(2, 3),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(5, 6),
    @(7, 7),
    @(6, 6),
    @(4, 4),
    @(6, 6),
    @(7, 7),
    @(4, 4),
    @(7, 7),
    @(4, 5),
    @(6, 6),
    @(8, 8),
    @(6, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
